# Slide 12, shape 2 ("2 Marcador de contenido") currently holds:
#   Paragraph 1: " Capacidad de " + "e" + "ditar" + " " + "s" + "ala"
#   Paragraph 2: "" (empty bullet)
#
# Target:
#   Paragraph 1: " Capacidad de editar sala"   (runs merged -> "editar" / "sala")
#   Paragraph 2 (new): " Exámen y reportes"
#   Paragraph 3: "" (empty bullet, unchanged)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

# --- Step 1: merge the split runs "e" + "ditar" -> "editar" ---
# Characters 15-20 (1-indexed) of paragraph 1 hold "e" followed by "ditar".
$eDitar = $tr.Characters(15, 6)
$eDitar.Text = "editar"

# --- Step 2: merge the split runs "s" + "ala" -> "sala" ---
# After the previous edit the text is unchanged in length, so "s" is still
# at position 22. Delete the lone "s" run, then type it back onto the
# front of the "ala" run so the surviving run keeps the "ala" run's
# formatting (which is what the target XML has: err="1" smtClean="0").
$sChar = $tr.Characters(22, 1)
$sChar.Text = ""
$ala = $tr.Characters(22, 3)
[void]$ala.InsertBefore("s")

# --- Step 3: insert a brand-new bullet paragraph after paragraph 1 ---
# Paragraph 2 is currently the trailing empty bullet line; push a new,
# empty paragraph in front of it (inherits the same bullet pPr) and then
# fill that new paragraph with the "Exámen y reportes" text.
$trailingEmpty = $tr.Paragraphs(2, 1)
[void]$trailingEmpty.InsertBefore("`r")

$newPara = $tr.Paragraphs(2, 1)
$newPara.Text = " Exámen y reportes"

# The insertion above leaves a stray empty run behind in what is now
# paragraph 3 (the original trailing empty bullet). Round-trip its text
# through a dummy value to flush that artifact back to a clean, run-less
# paragraph (matching the untouched original).
$oldTrailing = $tr.Paragraphs(3, 1)
$oldTrailing.Text = "x"
$oldTrailing2 = $tr.Paragraphs(3, 1)
$oldTrailing2.Text = ""
